$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (row 129), pushing the
# existing rows 129-209 down to 131-211. This mirrors the weekly refresh:
# the newest week's records are added at the front of the series and the
# oldest rows that fall past the previous bottom (208-209) reappear as the
# new bottom rows (210-211).
$ws.Rows.Item(129).Insert()
$ws.Rows.Item(129).Insert()

# New row 129 - "Primera" quality record for the new week (2022-01-21)
$ws.Cells.Item(129, 1).Value = 9
$ws.Cells.Item(129, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(129, 3).Value = "Metropolitana"
$ws.Cells.Item(129, 4).Value = 44582
$ws.Cells.Item(129, 5).Value = 13
$ws.Cells.Item(129, 6).Value = 100112017
$ws.Cells.Item(129, 7).Value = "Apio"
$ws.Cells.Item(129, 8).Value = "Americana (o)"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 70
$ws.Cells.Item(129, 11).Value = 6000
$ws.Cells.Item(129, 12).Value = 7000
$ws.Cells.Item(129, 13).Value = 6500
$ws.Cells.Item(129, 14).Value = "`$/docena de matas"
$ws.Cells.Item(129, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(129, 16).Value = 1083
$ws.Cells.Item(129, 17).Value = 6
$ws.Cells.Item(129, 18).Value = "Hortaliza"

# New row 130 - "Segunda" quality record for the same new week
$ws.Cells.Item(130, 1).Value = 9
$ws.Cells.Item(130, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(130, 3).Value = "Metropolitana"
$ws.Cells.Item(130, 4).Value = 44582
$ws.Cells.Item(130, 5).Value = 13
$ws.Cells.Item(130, 6).Value = 100112017
$ws.Cells.Item(130, 7).Value = "Apio"
$ws.Cells.Item(130, 8).Value = "Americana (o)"
$ws.Cells.Item(130, 9).Value = "Segunda"
$ws.Cells.Item(130, 10).Value = 34
$ws.Cells.Item(130, 11).Value = 5000
$ws.Cells.Item(130, 12).Value = 5000
$ws.Cells.Item(130, 13).Value = 5000
$ws.Cells.Item(130, 14).Value = "`$/docena de matas"
$ws.Cells.Item(130, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(130, 16).Value = 833
$ws.Cells.Item(130, 17).Value = 6
$ws.Cells.Item(130, 18).Value = "Hortaliza"
